$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 43: new "VS Code" entry ---
# Row 43 was the last (thick-bottom-bordered) row of the lookup table, and
# only had A43 (idx) + the shared H/I formulas + an empty J43 border cell.
# Row 27 is a fully-populated row using the same "section boundary" style
# (bottom border on A:G, J, center-aligned J, no border on K, bordered/
# centered L) - copy its formatting down onto row 43 first, then fill in
# the real content for the new "VS Code" / "vscodium" command.
$ws.Range("A27:L27").Copy()
$ws.Range("A43:L43").PasteSpecial(-4122)

$ws.Range("B43").Value = "VS Code"
$ws.Range("C43").Value = "vscodium"
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = 1
$ws.Range("J43").Value = "V"
$ws.Range("K43").Value = 86
$ws.Range("L43").Formula = "=DEC2HEX(K43)"

# --- Selection / scroll position on Sheet1 ---
$ws.Range("M29").Select()

# --- Workbook window placement (best effort; geometry of the host window) ---
$win = $excel.ActiveWindow
$win.Left = -20550
$win.Top = 5040
$win.Width = 16200
$win.Height = 28485
